$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# Row 2 ("electrical 4 u" reference): number cell becomes "4"
$tbl.Rows.Item(2).Cells.Item(1).Range.Text = "4"

# Row 3 (previously only held the _GoBack bookmark): number cell becomes "2",
# and the reference text cell gets the Jepson/Coleman/Igoe citation.
$tbl.Rows.Item(3).Cells.Item(1).Range.Text = "2"
$tbl.Rows.Item(3).Cells.Item(2).Range.Text = "B. Jepson, D. Coleman, and T. Igoe, 4. Introducing NDEF [Book]. Safari, 2017. [Online]. Available: https://www.safaribooksonline.com/library/view/beginning-nfc/9781449324094/ch04.html. Accessed: Feb. 2, 2017."

# Row 4 (previously fully empty): number cell becomes "3", text cell gets the
# NFC Forum / Business Wire citation.
$tbl.Rows.Item(4).Cells.Item(1).Range.Text = "3"
$tbl.Rows.Item(4).Cells.Item(2).Range.Text = "NFC, Forum. ""NFC Forum Technical Specifications Improve RF Communication and NFC Tag Interoperability with NFC Devices."" Business Wire (English) Dec. 0010: Regional Business News. Web. 1 Feb. 2017."

# Row 5 (previously fully empty): number cell becomes "5", text cell gets the
# NFC Forum operating-modes citation.
$tbl.Rows.Item(5).Cells.Item(1).Range.Text = "5"
$tbl.Rows.Item(5).Cells.Item(2).Range.Text = "N. Forum, A. rights reserved, A. M. services, Virtual, P. P. Terms, and C. Feedback, ""What are the operating modes of NFC devices? - NFC forum,"" NFC Forum, 2017. [Online]. Available: http://nfc-forum.org/resources/what-are-the-operating-modes-of-nfc-devices/. Accessed: Feb. 2, 2017."
